$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows with changed values
$ws.Range("D2").Value = "59.298.32"
$ws.Range("E2").Value = "  +2.57%  "
$ws.Range("D3").Value = "2.976.34"
$ws.Range("E3").Value = "  +1.46%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "561.77"
$ws.Range("E5").Value = "  +2.05%  "
$ws.Range("D6").Value = "137.61"
$ws.Range("E6").Value = "  +4.10%  "
$ws.Range("D8").Value = "0.519"
$ws.Range("E8").Value = "  +1.13%  "
$ws.Range("D9").Value = "2.962.36"
$ws.Range("E9").Value = "  +1.30%  "
$ws.Range("E10").Value = "  +4.00%  "
$ws.Range("D11").Value = "5.30"
$ws.Range("E11").Value = "  +11.20%  "
$ws.Range("D12").Value = "0.452"
$ws.Range("E12").Value = "  +1.70%  "
$ws.Range("D13").Value = "0.0000228"
$ws.Range("E13").Value = "  +3.78%  "
$ws.Range("D14").Value = "33.61"
$ws.Range("E14").Value = "  +2.55%  "
$ws.Range("E15").Value = "  -0.52%  "
$ws.Range("D16").Value = "3.471.06"
$ws.Range("E16").Value = "  +1.58%  "
$ws.Range("D17").Value = "7.13"
$ws.Range("E17").Value = "  +3.51%  "
$ws.Range("D18").Value = "2.976.94"
$ws.Range("E18").Value = "  +1.71%  "
$ws.Range("D19").Value = "59.349.84"
$ws.Range("D20").Value = "433.01"
$ws.Range("E20").Value = "  +4.21%  "
$ws.Range("D21").Value = "13.52"
$ws.Range("E21").Value = "  +1.62%  "
$ws.Range("D22").Value = "0.717"
$ws.Range("E22").Value = "  +3.31%  "
$ws.Range("D23").Value = "13.29"
$ws.Range("E23").Value = "  -1.21%  "
$ws.Range("D24").Value = "6.99"
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("D25").Value = "79.73"
$ws.Range("E25").Value = "  +0.55%  "
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("E27").Value = "  +10.17%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").Value = "2.53"
$ws.Range("E29").Value = "  +2.53%  "
$ws.Range("D30").Value = "7.72"
$ws.Range("E30").Value = "  +4.59%  "
$ws.Range("E31").Value = "  +8.52%  "
$ws.Range("D32").Value = "6.23"
$ws.Range("E32").Value = "  +4.74%  "
$ws.Range("D33").Value = "25.65"
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("E36").Value = "  +3.80%  "
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("D38").Value = "48.62"
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("D40").Value = "2.76"
$ws.Range("E40").Value = "  +5.85%  "
$ws.Range("D41").Value = "398.58"
$ws.Range("E41").Value = "  +6.23%  "
$ws.Range("E42").Value = "  +2.13%  "
$ws.Range("D43").Value = "2.747.95"
$ws.Range("E43").Value = "  +1.96%  "
$ws.Range("E44").Value = "  -2.73%  "
$ws.Range("D45").Value = "0.249"
$ws.Range("E45").Value = "  +6.02%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "122.67"
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("D48").Value = "34.16"
$ws.Range("E48").Value = "  +17.43%  "
$ws.Range("E49").Value = "  +1.67%  "
$ws.Range("E50").Value = "  +2.26%  "
$ws.Range("D51").Value = "23.26"
$ws.Range("E51").Value = "  +1.57%  "
# Row 34/35: PEPE and Mantle swap positions with updated values
$ws.Range("B34").Value = "Mantle"
$ws.Range("C34").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D34").Value = "0.987"
$ws.Range("E34").Value = "  +5.46%  "

$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").Value = "0.0₃0765"
$ws.Range("E35").Value = "  +10.08%  "
